{"js": "// Find the \"Built With\" bullet that currently reads just \"  - \" (the\n// trailing, still-empty bullet point) and append the \"Pixabay \u2013 Images.\"\n// citation to it, matching the author's addition. We reproduce the exact\n// run/proofErr structure (spell-check bookmarks around \"Pixabay\") used by\n// the sibling bullets (e.g. \"Coolors\", \"OpenAI-ChatGPT\") by inserting a\n// small OOXML fragment at the end of the paragraph's range.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the empty \"  - \" bullet paragraph (last item in the \"Built With\"\n// list that hasn't been filled in yet).\nlet target = null;\nfor (const p of paragraphs.items) {\n  if (p.text === \"  - \") {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find the \"  - \" placeholder bullet paragraph.');\n}\n\n// Collapse to the end of that paragraph's range so the new runs are\n// appended right after the existing \"  - \" run, inside the same paragraph.\nconst insertionPoint = target.getRange(Word.RangeLocation.end);\n\n// A minimal, well-formed FlatOPC package wrapping the <w:p> runs we want to\n// splice in. insertOoxml requires this pkg:package envelope.\nconst flatOpcXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p>\" +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  \"<w:r><w:t>Pixabay</w:t></w:r>\" +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> &#8211; Images.</w:t></w:r>' +\n  \"</w:p>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\ninsertionPoint.insertOoxml(flatOpcXml, Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Find the \"Built With\" bullet that is currently just \"  - \" (the last,\n# still-unfilled bullet in that list) and append the \"Pixabay \u2013 Images.\"\n# citation to it -- matching the run/proofErr (spell-check bookmark)\n# structure already used by its sibling bullets (\"Coolors\", \"OpenAI-\n# ChatGPT\", etc).\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph whose text (ignoring the trailing paragraph mark)\n# is exactly \"  - \".\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $paraText = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($paraText -eq \"  - \") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw 'Could not find the \"  - \" placeholder bullet paragraph.'\n}\n\n# Range covering just the paragraph's content, excluding its trailing\n# paragraph mark, so the new runs land inside the same paragraph.\n$r = $target.Range\n$r.End = $r.End - 1\n\n$dash = [char]0x2013\n$xml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n       '<w:proofErr w:type=\"spellStart\"/>' +\n       '<w:r><w:t>Pixabay</w:t></w:r>' +\n       '<w:proofErr w:type=\"spellEnd\"/>' +\n       '<w:r><w:t xml:space=\"preserve\"> ' + $dash + ' Images.</w:t></w:r>' +\n       '</w:p>'\n\n$r.InsertXML($xml)\n"}
